$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6886.625
$ws.Range("J40").Value = 7358.8
$ws.Range("L40").Value = 7358.8
$ws.Range("N40").Value = -7708.8

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null

$ws.Range("H100").Value = 1278.5714
$ws.Range("J100").Value = 1650
$ws.Range("L100").Value = 1650
$ws.Range("N100").Value = -2732

$ws.Range("H116").Value = 7999.6
$ws.Range("J116").Value = 8000
$ws.Range("L116").Value = 8000
$ws.Range("N116").Value = -14884

$ws.Range("H129").Value = 2088.6316
$ws.Range("I129").Value = 1021.375
$ws.Range("K129").Value = 3064.125
$ws.Range("M129").Value = 1935.875

$ws.Range("H132").Value = 6165
$ws.Range("I132").Value = 6518.9585
$ws.Range("K132").Value = 19556.8755
$ws.Range("M132").Value = -17026.8755

$ws.Range("H137").Value = 1067940.2
$ws.Range("I137").Value = 1725322.5
$ws.Range("J137").Value = 8824.388999999999
$ws.Range("K137").Value = 5175967.5
$ws.Range("L137").Value = 26473.167
$ws.Range("M137").Value = -5173417.5
$ws.Range("N137").Value = -31573.167

$ws.Range("H138").Value = 3351.3784
$ws.Range("J138").Value = 3707.1633
$ws.Range("L138").Value = 11121.4899
$ws.Range("N138").Value = -21401.4899

$ws.Range("H141").Value = 3661.75
$ws.Range("I141").Value = 3661.75
$ws.Range("K141").Value = 10985.25
$ws.Range("M141").Value = -5805.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 973.6786
$ws.Range("I2").Value = 890.2
$ws.Range("J2").Value = 1182.375
$ws.Range("K2").Value = 890.2
$ws.Range("L2").Value = 1182.375
$ws.Range("M2").Value = -777.2
$ws.Range("N2").Value = -1408.375

$ws.Range("H32").Value = 2666.6792
$ws.Range("I32").Value = 2666.6792
$ws.Range("K32").Value = 2666.6792
$ws.Range("M32").Value = -2379.6792

$ws.Range("H57").Value = 7525.5
$ws.Range("I57").Value = 7525.5
$ws.Range("K57").Value = 7525.5
$ws.Range("M57").Value = -7041.5

$ws.Range("H61").Value = 2778.0833
$ws.Range("I61").Value = 2279.76
$ws.Range("J61").Value = 3910.6365
$ws.Range("K61").Value = 2279.76
$ws.Range("L61").Value = 3910.6365
$ws.Range("M61").Value = -2067.76
$ws.Range("N61").Value = -4334.636500000001

$ws.Range("H62").Value = 49000
$ws.Range("J62").Value = 49000
$ws.Range("L62").Value = 49000
$ws.Range("N62").Value = -50248

$ws.Range("H65").Value = 49000
$ws.Range("J65").Value = 49000
$ws.Range("L65").Value = 147000
$ws.Range("N65").Value = -153240

$ws.Range("H74").Value = 127674.38
$ws.Range("I74").Value = 174505.53
$ws.Range("K74").Value = 174505.53
$ws.Range("M74").Value = -173631.53

$ws.Range("H77").Value = 127674.38
$ws.Range("I77").Value = 174505.53
$ws.Range("K77").Value = 872527.65
$ws.Range("M77").Value = -868159.65

$ws.Range("H102").Value = 10049.5
$ws.Range("I102").Value = 11779.6
$ws.Range("J102").Value = 7166
$ws.Range("K102").Value = 11779.6
$ws.Range("L102").Value = 7166
$ws.Range("M102").Value = -10157.6
$ws.Range("N102").Value = -10410

$ws.Range("H116").Value = 973.6786
$ws.Range("I116").Value = 890.2
$ws.Range("J116").Value = 1182.375
$ws.Range("K116").Value = 890.2
$ws.Range("L116").Value = 1182.375
$ws.Range("M116").Value = 1403.8
$ws.Range("N116").Value = -5770.375

$ws.Range("H122").Value = 2594.2
$ws.Range("I122").Value = 2636.1052
$ws.Range("J122").Value = 1798
$ws.Range("K122").Value = 7908.3156
$ws.Range("L122").Value = 5394
$ws.Range("M122").Value = -5458.3156
$ws.Range("N122").Value = -10294

$ws.Range("H132").Value = 2143.443
$ws.Range("I132").Value = 1830.5781
$ws.Range("J132").Value = 3478.3333
$ws.Range("K132").Value = 5491.7343
$ws.Range("L132").Value = 10434.9999
$ws.Range("M132").Value = -2961.7343
$ws.Range("N132").Value = -15494.9999

$ws.Range("H136").Value = 2778.0833
$ws.Range("I136").Value = 2279.76
$ws.Range("J136").Value = 3910.6365
$ws.Range("K136").Value = 6839.280000000001
$ws.Range("L136").Value = 11731.9095
$ws.Range("M136").Value = -4289.280000000001
$ws.Range("N136").Value = -16831.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 973.6786
$ws.Range("I3").Value = 890.2
$ws.Range("J3").Value = 1182.375
$ws.Range("K3").Value = 890.2
$ws.Range("L3").Value = 1182.375
$ws.Range("M3").Value = -776.2
$ws.Range("N3").Value = -1410.375

$ws.Range("H105").Value = 13002399
$ws.Range("I105").Value = 771067.4
$ws.Range("J105").Value = 35717730
$ws.Range("K105").Value = 771067.4
$ws.Range("L105").Value = 35717730
$ws.Range("M105").Value = -769320.4
$ws.Range("N105").Value = -35721224

$ws.Range("H134").Value = 3675.853
$ws.Range("I134").Value = 3442.68
$ws.Range("K134").Value = 10328.04
$ws.Range("M134").Value = -7793.039999999999

$ws.Range("H135").Value = 92251.336
$ws.Range("J135").Value = 92251.336
$ws.Range("L135").Value = 92251.336
$ws.Range("N135").Value = -102391.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3288.3076
$ws.Range("I31").Value = 2236.0488
$ws.Range("K31").Value = 2236.0488
$ws.Range("M31").Value = -1941.0488

$ws.Range("H34").Value = 3288.3076
$ws.Range("I34").Value = 2236.0488
$ws.Range("K34").Value = 2236.0488
$ws.Range("M34").Value = -2034.0488

$ws.Range("H52").Value = 93016.664
$ws.Range("J52").Value = 93016.664
$ws.Range("L52").Value = 93016.664
$ws.Range("N52").Value = -93604.664

$ws.Range("H58").Value = 2116.4443
$ws.Range("I58").Value = 1258
$ws.Range("J58").Value = 3833.3333
$ws.Range("K58").Value = 1258
$ws.Range("L58").Value = 3833.3333
$ws.Range("M58").Value = -1055
$ws.Range("N58").Value = -4239.3333

$ws.Range("H94").Value = 1903.8125
$ws.Range("I94").Value = 1898
$ws.Range("J94").Value = 1905.1538
$ws.Range("K94").Value = 1898
$ws.Range("L94").Value = 1905.1538
$ws.Range("M94").Value = -1447
$ws.Range("N94").Value = -2807.1538

$ws.Range("H132").Value = 20834728
$ws.Range("I132").Value = 27778738
$ws.Range("J132").Value = 2699.5
$ws.Range("K132").Value = 83336214
$ws.Range("L132").Value = 8098.5
$ws.Range("M132").Value = -83333684
$ws.Range("N132").Value = -13158.5

$ws.Range("H134").Value = 2414.9285
$ws.Range("I134").Value = 2130.762
$ws.Range("K134").Value = 6392.286
$ws.Range("M134").Value = -3857.286

$ws.Range("H136").Value = 2116.4443
$ws.Range("I136").Value = 1258
$ws.Range("J136").Value = 3833.3333
$ws.Range("K136").Value = 3774
$ws.Range("L136").Value = 11499.9999
$ws.Range("M136").Value = -1224
$ws.Range("N136").Value = -16599.9999

$ws.Range("H141").Value = 539106
$ws.Range("J141").Value = 539106
$ws.Range("L141").Value = 539106
$ws.Range("N141").Value = -549466

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 954.625
$ws.Range("J5").Value = 1983
$ws.Range("L5").Value = 5949
$ws.Range("N5").Value = -6173

$ws.Range("H8").Value = 719
$ws.Range("I8").Value = 719
$ws.Range("K8").Value = 2157
$ws.Range("M8").Value = -2018

$ws.Range("H135").Value = 954.625
$ws.Range("J135").Value = 1983
$ws.Range("L135").Value = 17847
$ws.Range("N135").Value = -22917

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8105.6
$ws.Range("I102").Value = 1357.1111
$ws.Range("J102").Value = 18228.334
$ws.Range("K102").Value = 1357.1111
$ws.Range("L102").Value = 18228.334
$ws.Range("M102").Value = 264.8888999999999
$ws.Range("N102").Value = -21472.334

$ws.Range("H107").Value = 6123.643
$ws.Range("I107").Value = 673.2
$ws.Range("K107").Value = 673.2
$ws.Range("M107").Value = 1246.8

$ws.Range("H122").Value = 4254.591
$ws.Range("I122").Value = 3002.2666
$ws.Range("J122").Value = 6938.143
$ws.Range("K122").Value = 9006.799800000001
$ws.Range("L122").Value = 20814.429
$ws.Range("M122").Value = -6556.799800000001
$ws.Range("N122").Value = -25714.429

$ws.Range("H123").Value = 78897.5
$ws.Range("J123").Value = 78897.5
$ws.Range("L123").Value = 78897.5
$ws.Range("N123").Value = -83797.5

$ws.Range("H132").Value = 2827.5757
$ws.Range("J132").Value = 6000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

$ws.Range("H134").Value = 68830.5
$ws.Range("J134").Value = 68830.5
$ws.Range("L134").Value = 206491.5
$ws.Range("N134").Value = -211561.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4657.2104
$ws.Range("I136").Value = 3655.4375
$ws.Range("K136").Value = 10966.3125
$ws.Range("M136").Value = -8416.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8336091.5
$ws.Range("I132").Value = 11113211
$ws.Range("K132").Value = 33339633
$ws.Range("M132").Value = -33337103

$ws.Range("H136").Value = 14404.473
$ws.Range("I136").Value = 14290.5
$ws.Range("J136").Value = 14700.8
$ws.Range("K136").Value = 42871.5
$ws.Range("L136").Value = 44102.39999999999
$ws.Range("M136").Value = -40321.5
$ws.Range("N136").Value = -49202.39999999999
